# Apply data-pipeline updates across Summary, Assets, and Liabilities sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B3").Value = "Hamad Al Qassimi"
$wsSummary.Range("B4").Value = 2447.82
$wsSummary.Range("B6").Value = 132333
$wsSummary.Range("B7").Value = 89857
$wsSummary.Range("B8").Value = 42476
$wsSummary.Range("B9").Value = 1.47

# ---------------------------------------------------------------------------
# Assets sheet
# ---------------------------------------------------------------------------
$wsAssets = $wb.Worksheets.Item("Assets")

$wsAssets.Range("B2").Value = "Mid-range Car"
$wsAssets.Range("C2").Value = 128655
$wsAssets.Range("C3").Value = 3678
$wsAssets.Range("C4").Value = 132333

# ---------------------------------------------------------------------------
# Liabilities sheet
# ---------------------------------------------------------------------------
$wsLiabilities = $wb.Worksheets.Item("Liabilities")

$wsLiabilities.Range("A2").Value = "Auto Loans"
$wsLiabilities.Range("B2").Value = "Vehicle Loan 1"
$wsLiabilities.Range("C2").Value = 77193
$wsLiabilities.Range("D2").Value = 919
$wsLiabilities.Range("E2").Value = 7

$wsLiabilities.Range("C3").Value = 12664
$wsLiabilities.Range("D3").Value = 633

$wsLiabilities.Range("C4").Value = 89857
